# "semana 18 de 2025" - add week-18 column (U), insert a new UPGD record row,
# and correct a couple of previously mis-shifted week-17 (T) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a brand-new row at 51 for UPGD 6600103414 / sub 01.
#    This pushes the old rows 51-53 (EPMSC PEREIRA, SANIDAD POLICIA
#    NACIONAL RISARALDA, BATALLON SAN MATEO) down to 52-54.
# ------------------------------------------------------------------
$ws.Rows(51).Insert()

$ws.Range("A51").Value = "'6600103414"
$ws.Range("B51").Value = "'01"

# ------------------------------------------------------------------
# 2) Add the new "18" (week 18) column header in U1, matching the
#    text style of the other week-number headers.
# ------------------------------------------------------------------
$ws.Range("U1").Value = "'18"

# ------------------------------------------------------------------
# 3) Fix up two pre-existing week-17 (T column) values.
# ------------------------------------------------------------------
$ws.Range("T32").Value = 23
$ws.Range("T40").Value = 169

# ------------------------------------------------------------------
# 4) Populate the new week-18 (U column) counts per UPGD row.
# ------------------------------------------------------------------
$weekU = @{
    2  = 29
    4  = 0
    5  = 4
    6  = 77
    7  = 32
    8  = 41
    10 = 2
    11 = 2
    12 = 2
    13 = 1
    17 = 2
    20 = 1
    22 = 2
    23 = 45
    26 = 124
    27 = 0
    28 = 15
    29 = 5
    30 = 5
    32 = 35
    33 = 4
    34 = 17
    35 = 83
    37 = 3
    38 = 27
    39 = 27
    40 = 153
    41 = 57
    42 = 302
    43 = 5
    44 = 88
    45 = 2
    46 = 0
    47 = 6
    48 = 2
    49 = 34
    50 = 0
    51 = 0
    52 = 2
    53 = 22
    54 = 33
}

foreach ($row in $weekU.Keys) {
    $ws.Range("U$row").Value = $weekU[$row]
}

Write-Output "week 18 column + new UPGD row applied"
